# Upload new version with timestamp
# Inserts a new inventory row for "LICID LOTION 30 ML" into the shortage
# report (alphabetically between "ITOMASH ..." and "LIMITLESS ..."),
# bumps the grand-total, and refreshes the generated-at timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 14 - this shifts every row from 14 downward
# (data rows, the totals row and the footer row) down by one, and copies
# formatting from the row above (row 13), matching the existing rows.
$ws.Rows("14:14").Insert()

# Fill in the new product row.
$ws.Cells.Item(14, 1).Value = 8              # A14 - sequence number
$ws.Cells.Item(14, 3).Value = "LICID LOTION 30 ML"   # C14 - item name
$ws.Cells.Item(14, 8).Value = "3:0"           # H14 - current balance
$ws.Cells.Item(14, 12).Value = 1              # L14 - order limit
$ws.Cells.Item(14, 14).Value = "40.00"        # N14 - price
$ws.Cells.Item(14, 16).Value = "40.0000"      # P14 - sell price
$ws.Cells.Item(14, 17).Value = "1:0"          # Q14 - transaction count

# Update the grand-total (previously row 21, now row 22 after the insert)
# to include the new item's price.
$ws.Cells.Item(22, 16).Value = 722.12

# Refresh the generated-at timestamp in the footer (previously row 22,
# now row 23 after the insert).
$ws.Cells.Item(23, 1).Value = "Tuesday, 9 September, 2025 11:17 AM"
